$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updates to existing F/G values in rows 699-761 (data corrections)
$ws.Range("F699").Value = 43482
$ws.Range("F700").Value = 43793
$ws.Range("F713").Value = 37317
$ws.Range("F719").Value = 44633
$ws.Range("F720").Value = 31239
$ws.Range("F723").Value = 22562
$ws.Range("F729").Value = 23229
$ws.Range("F730").Value = 19474
$ws.Range("F731").Value = 8612
$ws.Range("F732").Value = 11840
$ws.Range("F733").Value = 31640
$ws.Range("F734").Value = 23140
$ws.Range("F735").Value = 19310
$ws.Range("F739").Value = 8646
$ws.Range("F741").Value = 18855
$ws.Range("G741").Value = 1916
$ws.Range("F743").Value = 17974
$ws.Range("F744").Value = 14656
$ws.Range("F745").Value = 6144
$ws.Range("G745").Value = 924
$ws.Range("F746").Value = 7932
$ws.Range("G746").Value = 1233
$ws.Range("F747").Value = 22198
$ws.Range("G747").Value = 2350
$ws.Range("F748").Value = 16841
$ws.Range("G748").Value = 1524
$ws.Range("F749").Value = 14735
$ws.Range("G749").Value = 1468
$ws.Range("F750").Value = 14953
$ws.Range("G750").Value = 1340
$ws.Range("F751").Value = 12497
$ws.Range("G751").Value = 1375
$ws.Range("F752").Value = 4729
$ws.Range("G752").Value = 608
$ws.Range("F753").Value = 6719
$ws.Range("G753").Value = 931
$ws.Range("F754").Value = 20770
$ws.Range("G754").Value = 1908
$ws.Range("F755").Value = 13595
$ws.Range("G755").Value = 1271
$ws.Range("F756").Value = 13558
$ws.Range("G756").Value = 1047
$ws.Range("F757").Value = 13396
$ws.Range("G757").Value = 986
$ws.Range("F758").Value = 11113
$ws.Range("G758").Value = 920
$ws.Range("F759").Value = 3863
$ws.Range("G759").Value = 383
$ws.Range("F760").Value = 5077
$ws.Range("G760").Value = 549
$ws.Range("F761").Value = 16496
$ws.Range("G761").Value = 1233

# Fill in missing F/G values for row 762
$ws.Range("F762").Value = 10868
$ws.Range("G762").Value = 769

# Append new row 763 with data for 2022-04-06
$ws.Range("A763").Value = 44657
$ws.Range("B763").Value = 1741695
$ws.Range("C763").Value = 10952
$ws.Range("D763").Value = 4380
$ws.Range("E763").Value = 19500
$ws.Range("F763").Value = 6918
$ws.Range("G763").Value = 524
